$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 (LOGIN): just move the selection from B2 to B4.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B4").Select()

# ---------------------------------------------------------------------------
# Sheet 2 (3RD PARTY PROVIDER): replace the sample provider row with new
# glue-code test data, restyle it, and drop the now-unused trailing blank
# rows.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Select()

# New provider row values, written in the same order the original glue code
# appears to have used (this also controls shared-string allocation order).
$ws2.Range("H2").Value = "Justice Mahomed St & Steve Biko streets, Sunnyside, Pretoria, 0001, South Africa"
$ws2.Range("D2").Value = "xyz@gmail.com"
$ws2.Range("A2").Value = "xyz"
$ws2.Range("J2").Value = "Pretoria"
$ws2.Range("L2").Value = "South Africa"
$ws2.Range("K2").Value = "Gauteng"
$ws2.Range("N2").Value = "This is to create sample provider"
$ws2.Range("I2").Value = "959 Arcadia Street Hatfield, Pretoria, 0001, South Africa"
$ws2.Range("B2").Value = 9984738816
$ws2.Range("C2").Value = 8972736635
$ws2.Range("E2").Value = 7782937352
$ws2.Range("F2").Value = 9973263547
$ws2.Range("G2").Value = "Inactive"
$ws2.Range("M2").Value = 1111

# Remove the two now-empty trailing rows (old rows 3 & 4).
$ws2.Rows.Item(3).Delete()
$ws2.Rows.Item(3).Delete()

# --- Styling -----------------------------------------------------------
# Bulk "touched" style (applyAlignment flag, default alignment) across the
# header row and the bulk of the data row.
$ws2.Range("A1:N1").VerticalAlignment = -4107
$ws2.Range("B2:G2").VerticalAlignment = -4107
$ws2.Range("J2:L2").VerticalAlignment = -4107
$ws2.Range("N2").VerticalAlignment = -4107

# Vertically centered cells: provider name + street.
$ws2.Range("A2").VerticalAlignment = -4108
$ws2.Range("H2").VerticalAlignment = -4108

# Street2 cell: custom font (Arial, navy) + vertically centered.
$ws2.Range("I2").Font.Color = 5125145
$ws2.Range("I2").Font.Name = "Arial"
$ws2.Range("I2").VerticalAlignment = -4108

# Postal code column: custom number format "0;[Red]0".
$ws2.Range("M1:M2").NumberFormat = "0;[Red]0"

$ws2.Range("E2").Select()
